# Add two new test-case rows (batch_007 and batch_008) to Sheet1, mirroring
# the structure of the existing batch_001..batch_006 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Carry the formatting of row 6 down to the two new rows before filling
# them in, so the new cells pick up the same (text) number format and
# column-specific alignment used by the rest of the table (row 6 already
# has the "blank but formatted" E/G cells like the new rows need). Columns
# K/L are skipped since row 6 (like row 7) has no cells there.
$ws.Range("A6:J6").Copy() | Out-Null
$ws.Range("A8:J8").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:J9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("M6").Copy() | Out-Null
$ws.Range("M8").PasteSpecial(-4122) | Out-Null
$ws.Range("M9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 8: batch_007
$ws.Cells.Item(8, 1).Value = "batch_007"
$ws.Cells.Item(8, 2).Value = "y"
$ws.Cells.Item(8, 3).Value = "批量操作语句7执行"
$ws.Cells.Item(8, 4).Value = "batchsql"
$ws.Cells.Item(8, 6).Value = "batch07"
$ws.Cells.Item(8, 8).Value = "batch_sql_07"
$ws.Cells.Item(8, 9).Value = "select * from `$batch07"
$ws.Cells.Item(8, 10).Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_007.csv"
$ws.Cells.Item(8, 13).Value = "csv_containsAll"

# Row 9: batch_008
$ws.Cells.Item(9, 1).Value = "batch_008"
$ws.Cells.Item(9, 2).Value = "y"
$ws.Cells.Item(9, 3).Value = "批量操作语句8执行"
$ws.Cells.Item(9, 4).Value = "batchsql"
$ws.Cells.Item(9, 6).Value = "batch08"
$ws.Cells.Item(9, 8).Value = "batch_sql_08"
$ws.Cells.Item(9, 9).Value = "select * from `$batch08"
$ws.Cells.Item(9, 10).Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_008.csv"
$ws.Cells.Item(9, 13).Value = "csv_containsAll"

# Match the resulting active selection seen after the edit.
$ws.Range("J9").Select()
